$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest GitHub admin-log row (timestamp, action, org, team,
# actor, permission, pending-flag) into row 1 of the log sheet.
$ws.Range("A1").Value = "2025-07-23 11:39:30"
$ws.Range("B1").Value = "add-user"
$ws.Range("C1").Value = "new-organization97"
$ws.Range("D1").Value = "firstteam"
$ws.Range("F1").Value = "GokulJ17"
$ws.Range("G1").Value = "pull"
# Leading apostrophe forces Excel to store this as literal text "False"
# instead of auto-converting the look-alike word into a TRUE/FALSE boolean.
$ws.Range("I1").Value = "'False"
